# Atualização de bases das ligas, do dia: 23-02-2024 às 23:34
#
# Mexico Liga de Expansion.xlsx update:
#  1) Rows 7 and 8 had their match data (every column except the running
#     "id" in column A) swapped between each other.
#  2) Rows 262 and 263 had their match data (every column except the
#     running "id" in column A) swapped between each other.
#  3) Four brand-new match rows were appended at the bottom of the sheet
#     (the previous last row shifts down by one to make room for the
#     first new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) swap rows 7 / 8 (keep column A "id" sequence fixed) -----------------
$row7 = $ws.Range("B7:AC7")
$row8 = $ws.Range("B8:AC8")
$row7vals = $row7.Value()
$row8vals = $row8.Value()
$row7.Value = $row8vals
$row8.Value = $row7vals

# --- 2) swap rows 262 / 263 (keep column A "id" sequence fixed) ------------
$row262 = $ws.Range("B262:AC262")
$row263 = $ws.Range("B263:AC263")
$row262vals = $row262.Value()
$row263vals = $row263.Value()
$row262.Value = $row263vals
$row263.Value = $row262vals

# --- 3) append new rows at the bottom of the sheet -------------------------
# Copy the number formatting (bold/border id-style on col A, date style on
# col E) from the previous last row so the freshly written rows look the
# same as every other data row.
$ws.Range("A341").Copy()
$ws.Range("A342:A346").PasteSpecial(-4122)
$ws.Range("E341").Copy()
$ws.Range("E342:E346").PasteSpecial(-4122)

# New row 342
$ws.Range("A342").Value = 340
$ws.Range("B342").Value = 7641677
$ws.Range("C342").Value = "Mexico Liga de Expansion"
$ws.Range("D342").Value = "Mexico Liga de Expansion"
$ws.Range("E342").Value = 45345.00347222222
$ws.Range("F342").Value = "Club Atletico La Paz"
$ws.Range("G342").Value = "Dorados"
$ws.Range("H342").Value = 1
$ws.Range("I342").Value = 0
$ws.Range("J342").Value = "H"
$ws.Range("K342").Value = 1.6
$ws.Range("L342").Value = 4
$ws.Range("M342").Value = 4.75
$ws.Range("N342").Value = 1.7
$ws.Range("O342").Value = 4
$ws.Range("P342").Value = 4.75
$ws.Range("Q342").Value = -0.75
$ws.Range("R342").Value = 1.9
$ws.Range("S342").Value = 1.9
$ws.Range("T342").Value = 2.75
$ws.Range("U342").Value = 1.825
$ws.Range("V342").Value = 1.975
$ws.Range("W342").Value = 0.7
$ws.Range("X342").Value = -1
$ws.Range("Y342").Value = -1
$ws.Range("Z342").Value = 0.45
$ws.Range("AA342").Value = -0.5
$ws.Range("AB342").Value = -1
$ws.Range("AC342").Value = 0.9750000000000001

# Row 343 -- this is the sheet's former last row (old row 342), now
# shifted down to make room for the new row above; its data is unchanged.
$ws.Range("A343").Value = 341
$ws.Range("B343").Value = 7641680
$ws.Range("C343").Value = "Mexico Liga de Expansion"
$ws.Range("D343").Value = "Mexico Liga de Expansion"
$ws.Range("E343").Value = 45347.83680555555
$ws.Range("F343").Value = "Club Celaya"
$ws.Range("G343").Value = "Oaxaca"
$ws.Range("K343").Value = 1.333
$ws.Range("L343").Value = 4.75
$ws.Range("M343").Value = 7.5
$ws.Range("N343").Value = 1.363
$ws.Range("O343").Value = 5
$ws.Range("P343").Value = 8
$ws.Range("Q343").Value = -1.5
$ws.Range("R343").Value = 1.95
$ws.Range("S343").Value = 1.85
$ws.Range("T343").Value = 3
$ws.Range("U343").Value = 1.95
$ws.Range("V343").Value = 1.85
$ws.Range("W343").Value = 0
$ws.Range("X343").Value = 0
$ws.Range("Y343").Value = 0
$ws.Range("Z343").Value = 0
$ws.Range("AA343").Value = 0

# New row 344
$ws.Range("A344").Value = 342
$ws.Range("B344").Value = 7640646
$ws.Range("C344").Value = "Mexico Liga de Expansion"
$ws.Range("D344").Value = "Mexico Liga de Expansion"
$ws.Range("E344").Value = 45349.92013888889
$ws.Range("F344").Value = "Venados FC"
$ws.Range("G344").Value = "Atlante"
$ws.Range("K344").Value = 2.9
$ws.Range("L344").Value = 3.1
$ws.Range("M344").Value = 2.3
$ws.Range("N344").Value = 2.9
$ws.Range("O344").Value = 3.1
$ws.Range("P344").Value = 2.3
$ws.Range("Q344").Value = 0.25
$ws.Range("R344").Value = 1.775
$ws.Range("S344").Value = 2.025
$ws.Range("T344").Value = 2.25
$ws.Range("U344").Value = 2
$ws.Range("V344").Value = 1.8
$ws.Range("W344").Value = 0
$ws.Range("X344").Value = 0
$ws.Range("Y344").Value = 0
$ws.Range("Z344").Value = 0
$ws.Range("AA344").Value = 0

# New row 345
$ws.Range("A345").Value = 343
$ws.Range("B345").Value = 7641681
$ws.Range("C345").Value = "Mexico Liga de Expansion"
$ws.Range("D345").Value = "Mexico Liga de Expansion"
$ws.Range("E345").Value = 45350.00347222222
$ws.Range("F345").Value = "Universidad Guadalajara"
$ws.Range("G345").Value = "Correcaminos"
$ws.Range("K345").Value = 1.666
$ws.Range("L345").Value = 3.6
$ws.Range("M345").Value = 4.333
$ws.Range("N345").Value = 1.65
$ws.Range("O345").Value = 3.6
$ws.Range("P345").Value = 4.5
$ws.Range("Q345").Value = -0.75
$ws.Range("R345").Value = 1.85
$ws.Range("S345").Value = 1.95
$ws.Range("T345").Value = 2.5
$ws.Range("U345").Value = 1.8
$ws.Range("V345").Value = 2
$ws.Range("W345").Value = 0
$ws.Range("X345").Value = 0
$ws.Range("Y345").Value = 0
$ws.Range("Z345").Value = 0
$ws.Range("AA345").Value = 0

# New row 346
$ws.Range("A346").Value = 344
$ws.Range("B346").Value = 7641683
$ws.Range("C346").Value = "Mexico Liga de Expansion"
$ws.Range("D346").Value = "Mexico Liga de Expansion"
$ws.Range("E346").Value = 45351.00347222222
$ws.Range("F346").Value = "Dorados"
$ws.Range("G346").Value = "Tapatio"
$ws.Range("K346").Value = 2.9
$ws.Range("L346").Value = 3.4
$ws.Range("M346").Value = 2.15
$ws.Range("N346").Value = 2.9
$ws.Range("O346").Value = 3.4
$ws.Range("P346").Value = 2.15
$ws.Range("Q346").Value = 0.25
$ws.Range("R346").Value = 1.875
$ws.Range("S346").Value = 1.925
$ws.Range("T346").Value = 2.5
$ws.Range("U346").Value = 1.925
$ws.Range("V346").Value = 1.875
$ws.Range("W346").Value = 0
$ws.Range("X346").Value = 0
$ws.Range("Y346").Value = 0
$ws.Range("Z346").Value = 0
$ws.Range("AA346").Value = 0
